# Add two new users ("ids") to the Users table on Sheet1.
# Columns: A = Username, B = Password, C = ID
# Leading apostrophes force the numeric-looking password/ID values to be
# stored as text (matching the existing text-typed ID column), rather than
# being auto-converted to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "lala123#"
$ws.Range("B4").Value = "'1234"
$ws.Range("C4").Value = "'315783522"

$ws.Range("A5").Value = "lala123@"
$ws.Range("B5").Value = "matan123@$"
$ws.Range("C5").Value = "'315783522"
